$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column values are plain numeric-looking strings (e.g. "219.06").
# Force those specific cells to keep their existing Text representation so Excel
# does not silently convert them to numeric cell values.
$textCells = @('D5', 'D6', 'D8', 'D9', 'D11', 'D12', 'D15', 'D17', 'D19', 'D20', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D45', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.296.59'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '1.664.99'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.86%  '
$ws.Range('D5').Value = '219.06'
$ws.Range('D6').Value = '0.5333'
$ws.Range('E6').Value = '  +1.36%  '
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('D8').Value = '0.2646'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('D9').Value = '0.06413'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D11').Value = '0.07825'
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = '4.570'
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('D13').Value = '1.667.42'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').Value = '1.893.08'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').Value = '0.5523'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').Value = '0.0₅8207'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').Value = '65.64'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').Value = '4.698'
$ws.Range('E19').Value = '  +2.13%  '
$ws.Range('D20').Value = '193.47'
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').Value = '6.040'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').Value = '145.82'
$ws.Range('E24').Value = '  +2.77%  '
$ws.Range('D25').Value = '0.1234'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = '7.203'
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('D27').Value = '16.13'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = '1.483'
$ws.Range('E28').Value = '  +3.84%  '
$ws.Range('D29').Value = '0.05853'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('E30').Value = '  +0.50%  '
$ws.Range('D31').Value = '3.623'
$ws.Range('E31').Value = '  +3.05%  '
$ws.Range('D32').Value = '3.282'
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').Value = '1.609'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').Value = '0.9655'
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('D35').Value = '2.828'
$ws.Range('E35').Value = '  +1.61%  '
$ws.Range('D36').Value = '2.417'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').Value = '0.5807'
$ws.Range('E37').Value = '  +1.75%  '
$ws.Range('D38').Value = '0.01609'
$ws.Range('E38').Value = '  -0.82%  '
$ws.Range('D39').Value = '0.8667'
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('D40').Value = '5.881'
$ws.Range('E40').Value = '  +1.13%  '
$ws.Range('D41').Value = '1.051.19'
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('E42').Value = '  +0.76%  '
$ws.Range('E43').Value = '  +1.67%  '
$ws.Range('D44').Value = '1.803.63'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').Value = '57.74'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('E46').Value = '  -4.75%  '
$ws.Range('D47').Value = '1.014'
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('D48').Value = '0.4384'
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('D49').Value = '8.022'
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('D50').Value = '0.05166'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').Value = '1.416'
$ws.Range('E51').Value = '  -4.22%  '

# Restore the original (default) cell style so no stray formatting is introduced
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
